# Apply scheduled-runner updates to the Leviathan Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 822407.7
$ws.Range("I33").Value = 1015111.8
$ws.Range("J33").Value = 3415
$ws.Range("K33").Value = 1015111.8
$ws.Range("L33").Value = 3415
$ws.Range("M33").Value = -1014882.8
$ws.Range("N33").Value = -3873
$ws.Range("H40").Value = 2275.25
$ws.Range("I40").Value = 770.2
$ws.Range("J40").Value = 3350.2856
$ws.Range("K40").Value = 770.2
$ws.Range("L40").Value = 3350.2856
$ws.Range("M40").Value = -595.2
$ws.Range("N40").Value = -3700.2856
$ws.Range("H69").Value = 450
$ws.Range("I69").Value = 500
$ws.Range("J69").Value = 400
$ws.Range("K69").Value = 1500
$ws.Range("L69").Value = 1200
$ws.Range("M69").Value = -626
$ws.Range("N69").Value = -2948
$ws.Range("H72").Value = 450
$ws.Range("I72").Value = 500
$ws.Range("J72").Value = 400
$ws.Range("K72").Value = 4500
$ws.Range("L72").Value = 3600
$ws.Range("M72").Value = -132
$ws.Range("N72").Value = -12336
$ws.Range("H74").Value = 4071.7856
$ws.Range("I74").Value = 2301
$ws.Range("K74").Value = 2301
$ws.Range("M74").Value = -1365
$ws.Range("H77").Value = 4071.7856
$ws.Range("I77").Value = 2301
$ws.Range("K77").Value = 11505
$ws.Range("M77").Value = -6825
$ws.Range("H94").Value = 1339.5
$ws.Range("I94").Value = 1339.5
$ws.Range("K94").Value = 1339.5
$ws.Range("M94").Value = -888.5
$ws.Range("H96").Value = 200562.3
$ws.Range("I96").Value = 200562.3
$ws.Range("K96").Value = 601686.8999999999
$ws.Range("M96").Value = -600313.8999999999
$ws.Range("H97").Value = 1466.25
$ws.Range("J97").Value = 1418.6666
$ws.Range("L97").Value = 4255.9998
$ws.Range("N97").Value = -5247.9998
$ws.Range("H100").Value = 3150.4167
$ws.Range("I100").Value = 3130.5
$ws.Range("J100").Value = 3250
$ws.Range("K100").Value = 3130.5
$ws.Range("L100").Value = 3250
$ws.Range("M100").Value = -2589.5
$ws.Range("N100").Value = -4332
$ws.Range("H103").Value = 1475
$ws.Range("I103").Value = 900
$ws.Range("K103").Value = 2700
$ws.Range("M103").Value = -2114
$ws.Range("H112").Value = 2054.6
$ws.Range("J112").Value = 2073.5715
$ws.Range("L112").Value = 6220.7145
$ws.Range("N112").Value = -8436.7145
$ws.Range("H113").Value = 4752.5
$ws.Range("I113").Value = 3511.6667
$ws.Range("J113").Value = 5284.2856
$ws.Range("K113").Value = 3511.6667
$ws.Range("L113").Value = 5284.2856
$ws.Range("M113").Value = -257.6667000000002
$ws.Range("N113").Value = -11792.2856
$ws.Range("H138").Value = 3952.8333
$ws.Range("I138").Value = 2131.7144
$ws.Range("J138").Value = 5111.727
$ws.Range("K138").Value = 6395.1432
$ws.Range("L138").Value = 15335.181
$ws.Range("M138").Value = -1255.1432
$ws.Range("N138").Value = -25615.181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3034.3333
$ws.Range("I2").Value = 2496.7
$ws.Range("K2").Value = 2496.7
$ws.Range("M2").Value = -2383.7
$ws.Range("H11").Value = 4500
$ws.Range("I11").Value = 4500
$ws.Range("K11").Value = 4500
$ws.Range("M11").Value = -4356
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("H39").Value = 5062
$ws.Range("J39").Value = 12509
$ws.Range("L39").Value = 12509
$ws.Range("N39").Value = -13549
$ws.Range("H102").Value = 2865.9048
$ws.Range("I102").Value = 2364.7058
$ws.Range("K102").Value = 2364.7058
$ws.Range("M102").Value = -742.7058000000002
$ws.Range("H116").Value = 3034.3333
$ws.Range("I116").Value = 2496.7
$ws.Range("K116").Value = 2496.7
$ws.Range("M116").Value = -202.6999999999998
$ws.Range("H122").Value = 14093.059
$ws.Range("I122").Value = 16399.215
$ws.Range("J122").Value = 3331
$ws.Range("K122").Value = 49197.645
$ws.Range("L122").Value = 9993
$ws.Range("M122").Value = -46747.645
$ws.Range("N122").Value = -14893

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3034.3333
$ws.Range("I3").Value = 2496.7
$ws.Range("K3").Value = 2496.7
$ws.Range("M3").Value = -2382.7
$ws.Range("H35").Value = 30000
$ws.Range("J35").Value = 30000
$ws.Range("L35").Value = 30000
$ws.Range("N35").Value = -30620
$ws.Range("H105").Value = 6252170
$ws.Range("J105").Value = 1999
$ws.Range("L105").Value = 1999
$ws.Range("N105").Value = -5493
$ws.Range("H132").Value = 178850
$ws.Range("J132").Value = 178850
$ws.Range("L132").Value = 178850
$ws.Range("N132").Value = -188970

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5784.2383
$ws.Range("I58").Value = 7574.2
$ws.Range("J58").Value = 1309.3334
$ws.Range("K58").Value = 7574.2
$ws.Range("L58").Value = 1309.3334
$ws.Range("M58").Value = -7371.2
$ws.Range("N58").Value = -1715.3334
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H103").Value = 19027.666
$ws.Range("I103").Value = 19027.666
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 19027.666
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -17855.666
$ws.Range("N103").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 5784.2383
$ws.Range("I136").Value = 7574.2
$ws.Range("J136").Value = 1309.3334
$ws.Range("K136").Value = 22722.6
$ws.Range("L136").Value = 3928.0002
$ws.Range("M136").Value = -20172.6
$ws.Range("N136").Value = -9028.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 62999
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H50").Value = 62999
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H122").Value = 4637.9165
$ws.Range("I122").Value = 4426.375
$ws.Range("K122").Value = 13279.125
$ws.Range("M122").Value = -10829.125
$ws.Range("H132").Value = 2064.7
$ws.Range("J132").Value = 2181.25
$ws.Range("L132").Value = 6543.75
$ws.Range("N132").Value = -11603.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()
$ws.Range("H46").Value = 55352.125
$ws.Range("I46").Value = 72302.836
$ws.Range("J46").Value = 4500
$ws.Range("K46").Value = 72302.836
$ws.Range("L46").Value = 4500
$ws.Range("M46").Value = -72114.836
$ws.Range("N46").Value = -4876
$ws.Range("H47").Value = 15000
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H52").Value = 15000
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H61").Value = 22248032
$ws.Range("I61").Value = 27781674
$ws.Range("J61").Value = 113464.336
$ws.Range("K61").Value = 27781674
$ws.Range("L61").Value = 113464.336
$ws.Range("M61").Value = -27781472
$ws.Range("N61").Value = -113868.336
$ws.Range("H113").Value = 22248032
$ws.Range("I113").Value = 27781674
$ws.Range("J113").Value = 113464.336
$ws.Range("K113").Value = 27781674
$ws.Range("L113").Value = 113464.336
$ws.Range("M113").Value = -27779504
$ws.Range("N113").Value = -117804.336
$ws.Range("H122").Value = 4685.5713
$ws.Range("J122").Value = 4635
$ws.Range("L122").Value = 13905
$ws.Range("N122").Value = -18805

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 116570.79
$ws.Range("I62").Value = 4571
$ws.Range("K62").Value = 4571
$ws.Range("M62").Value = -3947
$ws.Range("H65").Value = 116570.79
$ws.Range("I65").Value = 4571
$ws.Range("K65").Value = 22855
$ws.Range("M65").Value = -19735
$ws.Range("H100").Value = 3176.2
$ws.Range("I100").Value = 4404.375
$ws.Range("J100").Value = 1772.5714
$ws.Range("K100").Value = 8808.75
$ws.Range("L100").Value = 3545.1428
$ws.Range("M100").Value = -8267.75
